$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.062550666666667
$ws.Range("H2").Value = 9.187652
$ws.Range("I2").Value = 0.06849600470812313
$ws.Range("J2").Value = 0.06849600470812313
$ws.Range("M2").Value = 119.0164006666667
$ws.Range("N2").Value = 357.049202
$ws.Range("O2").Value = 0.9176278005170622
$ws.Range("P2").Value = 0.9176278005170622
$ws.Range("Q2").Value = 364.4937572059671
$ws.Range("R2").Value = 3280.443814853704
$ws.Range("S2").Value = 0.06285383814452136
$ws.Range("T2").Value = 0.06285383814452136
$ws.Range("G3").Value = 3.062550666666667
$ws.Range("H3").Value = 9.187652
$ws.Range("I3").Value = 0.06849600470812313
$ws.Range("J3").Value = 0.06849600470812313
$ws.Range("O3").Value = 0.001755838010330732
$ws.Range("P3").Value = 0.001755838010330731
$ws.Range("Q3").Value = 0.6974418092715555
$ws.Range("R3").Value = 6.276976283444
$ws.Range("S3").Value = 0.0001202678886223153
$ws.Range("T3").Value = 0.0001202678886223153
$ws.Range("G4").Value = 3.062550666666667
$ws.Range("H4").Value = 9.187652
$ws.Range("I4").Value = 0.06849600470812313
$ws.Range("J4").Value = 0.06849600470812313
$ws.Range("M4").Value = 7.816301333333333
$ws.Range("N4").Value = 23.448904
$ws.Range("O4").Value = 0.06026442877207647
$ws.Range("P4").Value = 0.06026442877207646
$ws.Range("Q4").Value = 23.93781885926755
$ws.Range("R4").Value = 215.440369733408
$ws.Range("S4").Value = 0.004127872596904501
$ws.Range("T4").Value = 0.0041278725969045
$ws.Range("G5").Value = 3.062550666666667
$ws.Range("H5").Value = 9.187652
$ws.Range("I5").Value = 0.06849600470812313
$ws.Range("J5").Value = 0.06849600470812313
$ws.Range("M5").Value = 0.105045
$ws.Range("N5").Value = 0.315135
$ws.Range("O5").Value = 0.000809906968832672
$ws.Range("P5").Value = 0.000809906968832672
$ws.Range("Q5").Value = 0.32170563478
$ws.Range("R5").Value = 2.89535071302
$ws.Range("S5").Value = 0.00005547539155030443
$ws.Range("T5").Value = 0.00005547539155030443
$ws.Range("G6").Value = 3.062550666666667
$ws.Range("H6").Value = 9.187652
$ws.Range("I6").Value = 0.06849600470812313
$ws.Range("J6").Value = 0.06849600470812313
$ws.Range("M6").Value = 2.534602333333333
$ws.Range("N6").Value = 7.603807
$ws.Range("O6").Value = 0.01954202573169801
$ws.Range("P6").Value = 0.01954202573169801
$ws.Range("Q6").Value = 7.762348065684889
$ws.Range("R6").Value = 69.86113259116399
$ws.Range("S6").Value = 0.00133855068652465
$ws.Range("T6").Value = 0.00133855068652465
$ws.Range("I7").Value = 0.396815038797359
$ws.Range("J7").Value = 0.396815038797359
$ws.Range("M7").Value = 119.0164006666667
$ws.Range("N7").Value = 357.049202
$ws.Range("O7").Value = 0.9176278005170622
$ws.Range("P7").Value = 0.9176278005170622
$ws.Range("Q7").Value = 2111.606436366764
$ws.Range("R7").Value = 19004.45792730088
$ws.Range("S7").Value = 0.3641285112637133
$ws.Range("T7").Value = 0.3641285112637133
$ws.Range("I8").Value = 0.396815038797359
$ws.Range("J8").Value = 0.396815038797359
$ws.Range("O8").Value = 0.001755838010330732
$ws.Range("P8").Value = 0.001755838010330731
$ws.Range("S8").Value = 0.0006967429281912669
$ws.Range("T8").Value = 0.0006967429281912668
$ws.Range("I9").Value = 0.396815038797359
$ws.Range("J9").Value = 0.396815038797359
$ws.Range("M9").Value = 7.816301333333333
$ws.Range("N9").Value = 23.448904
$ws.Range("O9").Value = 0.06026442877207647
$ws.Range("P9").Value = 0.06026442877207646
$ws.Range("Q9").Value = 138.6779646468622
$ws.Range("R9").Value = 1248.10168182176
$ws.Range("S9").Value = 0.0239138316412922
$ws.Range("T9").Value = 0.0239138316412922
$ws.Range("I10").Value = 0.396815038797359
$ws.Range("J10").Value = 0.396815038797359
$ws.Range("M10").Value = 0.105045
$ws.Range("N10").Value = 0.315135
$ws.Range("O10").Value = 0.000809906968832672
$ws.Range("P10").Value = 0.000809906968832672
$ws.Range("Q10").Value = 1.8637237966
$ws.Range("R10").Value = 16.7735141694
$ws.Range("S10").Value = 0.0003213832652595882
$ws.Range("T10").Value = 0.0003213832652595882
$ws.Range("I11").Value = 0.396815038797359
$ws.Range("J11").Value = 0.396815038797359
$ws.Range("M11").Value = 2.534602333333333
$ws.Range("N11").Value = 7.603807
$ws.Range("O11").Value = 0.01954202573169801
$ws.Range("P11").Value = 0.01954202573169801
$ws.Range("Q11").Value = 44.96928633967556
$ws.Range("R11").Value = 404.72357705708
$ws.Range("S11").Value = 0.007754569698902735
$ws.Range("T11").Value = 0.007754569698902735
$ws.Range("G12").Value = 13.27534766666667
$ws.Range("H12").Value = 39.826043
$ws.Range("I12").Value = 0.2969120759943797
$ws.Range("J12").Value = 0.2969120759943796
$ws.Range("M12").Value = 119.0164006666667
$ws.Range("N12").Value = 357.049202
$ws.Range("O12").Value = 0.9176278005170622
$ws.Range("P12").Value = 0.9176278005170622
$ws.Range("Q12").Value = 1579.984096885298
$ws.Range("R12").Value = 14219.85687196768
$ws.Range("S12").Value = 0.2724547752416774
$ws.Range("T12").Value = 0.2724547752416774
$ws.Range("G13").Value = 13.27534766666667
$ws.Range("H13").Value = 39.826043
$ws.Range("I13").Value = 0.2969120759943797
$ws.Range("J13").Value = 0.2969120759943796
$ws.Range("O13").Value = 0.001755838010330732
$ws.Range("P13").Value = 0.001755838010330731
$ws.Range("Q13").Value = 3.023225899941222
$ws.Range("R13").Value = 27.209033099471
$ws.Range("S13").Value = 0.0005213295087571386
$ws.Range("T13").Value = 0.0005213295087571384
$ws.Range("G14").Value = 13.27534766666667
$ws.Range("H14").Value = 39.826043
$ws.Range("I14").Value = 0.2969120759943797
$ws.Range("J14").Value = 0.2969120759943796
$ws.Range("M14").Value = 7.816301333333333
$ws.Range("N14").Value = 23.448904
$ws.Range("O14").Value = 0.06026442877207647
$ws.Range("P14").Value = 0.06026442877207646
$ws.Range("Q14").Value = 103.7641176674302
$ws.Range("R14").Value = 933.8770590068719
$ws.Range("S14").Value = 0.01789323665533265
$ws.Range("T14").Value = 0.01789323665533265
$ws.Range("G15").Value = 13.27534766666667
$ws.Range("H15").Value = 39.826043
$ws.Range("I15").Value = 0.2969120759943797
$ws.Range("J15").Value = 0.2969120759943796
$ws.Range("M15").Value = 0.105045
$ws.Range("N15").Value = 0.315135
$ws.Range("O15").Value = 0.000809906968832672
$ws.Range("P15").Value = 0.000809906968832672
$ws.Range("Q15").Value = 1.394508895645
$ws.Range("R15").Value = 12.550580060805
$ws.Range("S15").Value = 0.000240471159478424
$ws.Range("T15").Value = 0.000240471159478424
$ws.Range("G16").Value = 13.27534766666667
$ws.Range("H16").Value = 39.826043
$ws.Range("I16").Value = 0.2969120759943797
$ws.Range("J16").Value = 0.2969120759943796
$ws.Range("M16").Value = 2.534602333333333
$ws.Range("N16").Value = 7.603807
$ws.Range("O16").Value = 0.01954202573169801
$ws.Range("P16").Value = 0.01954202573169801
$ws.Range("Q16").Value = 33.64772717174456
$ws.Range("R16").Value = 302.829544545701
$ws.Range("S16").Value = 0.005802263429134044
$ws.Range("T16").Value = 0.005802263429134043
$ws.Range("G17").Value = 3.455866
$ws.Range("H17").Value = 10.367598
$ws.Range("I17").Value = 0.07729276657626213
$ws.Range("J17").Value = 0.07729276657626213
$ws.Range("M17").Value = 119.0164006666667
$ws.Range("N17").Value = 357.049202
$ws.Range("O17").Value = 0.9176278005170622
$ws.Range("P17").Value = 0.9176278005170622
$ws.Range("Q17").Value = 411.3047325063106
$ws.Range("R17").Value = 3701.742592556795
$ws.Range("S17").Value = 0.07092599138925412
$ws.Range("T17").Value = 0.07092599138925412
$ws.Range("G18").Value = 3.455866
$ws.Range("H18").Value = 10.367598
$ws.Range("I18").Value = 0.07729276657626213
$ws.Range("J18").Value = 0.07729276657626213
$ws.Range("O18").Value = 0.001755838010330732
$ws.Range("P18").Value = 0.001755838010330731
$ws.Range("Q18").Value = 0.7870124278673334
$ws.Range("R18").Value = 7.083111850806
$ws.Range("S18").Value = 0.0001357135774782218
$ws.Range("T18").Value = 0.0001357135774782217
$ws.Range("G19").Value = 3.455866
$ws.Range("H19").Value = 10.367598
$ws.Range("I19").Value = 0.07729276657626213
$ws.Range("J19").Value = 0.07729276657626213
$ws.Range("M19").Value = 7.816301333333333
$ws.Range("N19").Value = 23.448904
$ws.Range("O19").Value = 0.06026442877207647
$ws.Range("P19").Value = 0.06026442877207646
$ws.Range("Q19").Value = 27.01209002362133
$ws.Range("R19").Value = 243.108810212592
$ws.Range("S19").Value = 0.004658004425931882
$ws.Range("T19").Value = 0.004658004425931881
$ws.Range("G20").Value = 3.455866
$ws.Range("H20").Value = 10.367598
$ws.Range("I20").Value = 0.07729276657626213
$ws.Range("J20").Value = 0.07729276657626213
$ws.Range("M20").Value = 0.105045
$ws.Range("N20").Value = 0.315135
$ws.Range("O20").Value = 0.000809906968832672
$ws.Range("P20").Value = 0.000809906968832672
$ws.Range("Q20").Value = 0.36302144397
$ws.Range("R20").Value = 3.26719299573
$ws.Range("S20").Value = 0.00006259995029047172
$ws.Range("T20").Value = 0.00006259995029047172
$ws.Range("G21").Value = 3.455866
$ws.Range("H21").Value = 10.367598
$ws.Range("I21").Value = 0.07729276657626213
$ws.Range("J21").Value = 0.07729276657626213
$ws.Range("M21").Value = 2.534602333333333
$ws.Range("N21").Value = 7.603807
$ws.Range("O21").Value = 0.01954202573169801
$ws.Range("P21").Value = 0.01954202573169801
$ws.Range("Q21").Value = 8.759246027287332
$ws.Range("R21").Value = 78.83321424558599
$ws.Range("S21").Value = 0.001510457233307443
$ws.Range("T21").Value = 0.001510457233307443
$ws.Range("G22").Value = 7.175465666666668
$ws.Range("H22").Value = 21.526397
$ws.Range("I22").Value = 0.1604841139238761
$ws.Range("J22").Value = 0.1604841139238761
$ws.Range("M22").Value = 119.0164006666667
$ws.Range("N22").Value = 357.049202
$ws.Range("O22").Value = 0.9176278005170622
$ws.Range("P22").Value = 0.9176278005170622
$ws.Range("Q22").Value = 853.9980967539104
$ws.Range("R22").Value = 7685.982870785195
$ws.Range("S22").Value = 0.147264684477896
$ws.Range("T22").Value = 0.147264684477896
$ws.Range("G23").Value = 7.175465666666668
$ws.Range("H23").Value = 21.526397
$ws.Range("I23").Value = 0.1604841139238761
$ws.Range("J23").Value = 0.1604841139238761
$ws.Range("O23").Value = 0.001755838010330732
$ws.Range("P23").Value = 0.001755838010330731
$ws.Range("Q23").Value = 1.634085539023223
$ws.Range("R23").Value = 14.706769851209
$ws.Range("S23").Value = 0.000281784107281789
$ws.Range("T23").Value = 0.000281784107281789
$ws.Range("G24").Value = 7.175465666666668
$ws.Range("H24").Value = 21.526397
$ws.Range("I24").Value = 0.1604841139238761
$ws.Range("J24").Value = 0.1604841139238761
$ws.Range("M24").Value = 7.816301333333333
$ws.Range("N24").Value = 23.448904
$ws.Range("O24").Value = 0.06026442877207647
$ws.Range("P24").Value = 0.06026442877207646
$ws.Range("Q24").Value = 56.08560185765423
$ws.Range("R24").Value = 504.770416718888
$ws.Range("S24").Value = 0.009671483452615235
$ws.Range("T24").Value = 0.009671483452615233
$ws.Range("G25").Value = 7.175465666666668
$ws.Range("H25").Value = 21.526397
$ws.Range("I25").Value = 0.1604841139238761
$ws.Range("J25").Value = 0.1604841139238761
$ws.Range("M25").Value = 0.105045
$ws.Range("N25").Value = 0.315135
$ws.Range("O25").Value = 0.000809906968832672
$ws.Range("P25").Value = 0.000809906968832672
$ws.Range("Q25").Value = 0.7537467909550001
$ws.Range("R25").Value = 6.783721118595001
$ws.Range("S25").Value = 0.0001299772022538837
$ws.Range("T25").Value = 0.0001299772022538837
$ws.Range("G26").Value = 7.175465666666668
$ws.Range("H26").Value = 21.526397
$ws.Range("I26").Value = 0.1604841139238761
$ws.Range("J26").Value = 0.1604841139238761
$ws.Range("M26").Value = 2.534602333333333
$ws.Range("N26").Value = 7.603807
$ws.Range("O26").Value = 0.01954202573169801
$ws.Range("P26").Value = 0.01954202573169801
$ws.Range("Q26").Value = 18.18695202148656
$ws.Range("R26").Value = 163.682568193379
$ws.Range("S26").Value = 0.003136184683829142
$ws.Range("T26").Value = 0.003136184683829142
